# TC02_C3DC_phs001437_SexAtBirth-Female: fix the Treatment Agent query
# (drop the redundant CONCAT() wrapper around REPLACE()) and restyle the
# cell to match the rest of the "Added C3DC phs001437 to Smoke suite" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B5")

$oldFragment = "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"","
$newFragment = "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","

$cell.Value2 = $cell.Value2.Replace($oldFragment, $newFragment)

# The query cell also picked up an explicit 11pt font (the other TabQuery
# cells in column B stay at their original 12pt).
$cell.Font.Size = 11
